$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (Standard Deviation) updated values
$ws.Range("B3").Value = 0.001753054283238868
$ws.Range("C3").Value = 0.001042925767129162
$ws.Range("D3").Value = 0.008212652505215173
$ws.Range("E3").Value = 0.01755742829292724
$ws.Range("F3").Value = 0.007386946580121634
$ws.Range("G3").Value = 0.04066989668766503

# Row 4 (Maximum) updated value
$ws.Range("B4").Value = 0.005324698759889668
